$wb = $excel.ActiveWorkbook

# Overview sheet: update "Latest HO Xliff Generate Date" timestamp
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-23 21:09:29"

# zh-cn sheet: update "Correspond Handoff Datetime" and "Correspond Handback DateTime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-23 21:09:24"
$wsZhCn.Range("K2").Value = "2016-08-23 21:09:41"

# de-de sheet: "Correspond Handoff Datetime" (H2) shares the same underlying
# string as Overview!G2 (both were "2016-08-23 21:08:43"), so it moves to the
# same new value; "Correspond Handback DateTime" (K2) gets its own new value.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-23 21:09:29"
$wsDeDe.Range("K2").Value = "2016-08-23 21:09:49"
